$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '245.26'
Set-TextValue 'D3' '25.12'
Set-TextValue 'D4' '5.028'
Set-TextValue 'D5' '0.05618'
Set-TextValue 'D6' '6.578'
Set-TextValue 'D7' '3.011'
Set-TextValue 'D9' '0.8364'
Set-TextValue 'D10' '0.1337'
Set-TextValue 'D13' '0.09405'
Set-TextValue 'D14' '0.001531'
Set-TextValue 'D15' '0.0005946'
Set-TextValue 'E15' '14OneONEWorstin24h'
Set-TextValue 'D16' '0.006091'
Set-TextValue 'D17' '3.499'
Set-TextValue 'D18' '2.092'
Set-TextValue 'D20' '0.03252'
Set-TextValue 'D21' '0.1292'
Set-TextValue 'D22' '3.745'
Set-TextValue 'D23' '0.04664'
Set-TextValue 'D24' '0.1369'
Set-TextValue 'D25' '0.001242'
Set-TextValue 'D27' '0.00009689'
Set-TextValue 'D28' '0.0001939'
Set-TextValue 'B41' 'KickToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.006223'
Set-TextValue 'E41' '40KickTokenKICK'
Set-TextValue 'B42' 'BKEXToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1056'
Set-TextValue 'E42' '41BKEXTokenBKK'
Set-TextValue 'B43' 'CEJI'
Set-TextValue 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue 'D43' '0.002732'
Set-TextValue 'E43' '42CEJICEJI'
Set-TextValue 'D44' '0.008164'
Set-TextValue 'D45' '0.00005293'
Set-TextValue 'D47' '0.2259'
Set-TextValue 'D48' '0.002019'
